# Scheduled-runner style refresh of market-price-derived columns (H:N)
# across the per-job profit sheets. Only numeric value cells change;
# row/column layout, tables and text columns are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 257.5
$ws.Range("I18").Value = 257.5
$ws.Range("K18").Value = 257.5
$ws.Range("M18").Value = 26.5

$ws.Range("H19").Value = 2654.3
$ws.Range("I19").Value = 3122.75
$ws.Range("J19").Value = 780.5
$ws.Range("K19").Value = 3122.75
$ws.Range("L19").Value = 780.5
$ws.Range("M19").Value = -2947.75
$ws.Range("N19").Value = -1130.5

$ws.Range("H80").Value = 836531.25
$ws.Range("I80").Value = 1779.8334
$ws.Range("K80").Value = 5339.5002
$ws.Range("M80").Value = -4341.5002

$ws.Range("H83").Value = 836531.25
$ws.Range("I83").Value = 1779.8334
$ws.Range("K83").Value = 16018.5006
$ws.Range("M83").Value = -11026.5006

$ws.Range("H87").Value = 81662
$ws.Range("J87").Value = 84540.37
$ws.Range("L87").Value = 84540.37
$ws.Range("N87").Value = -87036.37

$ws.Range("H90").Value = 81662
$ws.Range("J90").Value = 84540.37
$ws.Range("L90").Value = 253621.11
$ws.Range("N90").Value = -266101.11

$ws.Range("H101").Value = 854.5833
$ws.Range("I101").Value = 833.1111
$ws.Range("J101").Value = 919
$ws.Range("K101").Value = 2499.3333
$ws.Range("L101").Value = 2757
$ws.Range("M101").Value = -877.3332999999998
$ws.Range("N101").Value = -6001

$ws.Range("H111").Value = 2437.875
$ws.Range("I111").Value = 2107.5715
$ws.Range("K111").Value = 6322.7145
$ws.Range("M111").Value = -3255.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3282.923
$ws.Range("I63").Value = 2152
$ws.Range("K63").Value = 2152
$ws.Range("M63").Value = -1466

$ws.Range("H66").Value = 3282.923
$ws.Range("I66").Value = 2152
$ws.Range("K66").Value = 10760
$ws.Range("M66").Value = -7328

$ws.Range("H122").Value = 1824
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 1848
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 5544
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -10444

$ws.Range("H139").Value = 73571
$ws.Range("J139").Value = 73571
$ws.Range("L139").Value = 73571
$ws.Range("N139").Value = -83851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 480119.1
$ws.Range("I86").Value = 1002756.8
$ws.Range("K86").Value = 1002756.8
$ws.Range("M86").Value = -1001633.8

$ws.Range("H89").Value = 480119.1
$ws.Range("I89").Value = 1002756.8
$ws.Range("K89").Value = 5013784
$ws.Range("M89").Value = -5008168

$ws.Range("H99").Value = 8641.846
$ws.Range("I99").Value = 8528.75
$ws.Range("K99").Value = 8528.75
$ws.Range("M99").Value = -7030.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2766.0908
$ws.Range("I31").Value = 1694.7273
$ws.Range("K31").Value = 1694.7273
$ws.Range("M31").Value = -1399.7273

$ws.Range("H34").Value = 2766.0908
$ws.Range("I34").Value = 1694.7273
$ws.Range("K34").Value = 1694.7273
$ws.Range("M34").Value = -1492.7273

$ws.Range("H51").Value = 32660.223
$ws.Range("J51").Value = 32660.223
$ws.Range("L51").Value = 32660.223
$ws.Range("N51").Value = -34132.223

$ws.Range("H60").Value = 32396.5
$ws.Range("I60").Value = 13498.25
$ws.Range("J60").Value = 44995.332
$ws.Range("K60").Value = 13498.25
$ws.Range("L60").Value = 44995.332
$ws.Range("M60").Value = -12987.25
$ws.Range("N60").Value = -46017.332

$ws.Range("H61").Value = 32660.223
$ws.Range("J61").Value = 32660.223
$ws.Range("L61").Value = 32660.223
$ws.Range("N61").Value = -33356.223

$ws.Range("H107").Value = 834.9032
$ws.Range("I107").Value = 822.88
$ws.Range("K107").Value = 822.88
$ws.Range("M107").Value = 1097.12

$ws.Range("H122").Value = 2027.5
$ws.Range("I122").Value = 1019.1111
$ws.Range("J122").Value = 5052.6665
$ws.Range("K122").Value = 3057.3333
$ws.Range("L122").Value = 15157.9995
$ws.Range("M122").Value = -607.3332999999998
$ws.Range("N122").Value = -20057.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2040.1428
$ws.Range("I107").Value = 1926.75
$ws.Range("J107").Value = 2191.3333
$ws.Range("K107").Value = 5780.25
$ws.Range("L107").Value = 6573.999899999999
$ws.Range("M107").Value = -3860.25
$ws.Range("N107").Value = -10413.9999

$ws.Range("H122").Value = 20000856
$ws.Range("I122").Value = 749
$ws.Range("J122").Value = 25000882
$ws.Range("K122").Value = 6741
$ws.Range("L122").Value = 225007938
$ws.Range("M122").Value = -4291
$ws.Range("N122").Value = -225012838

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1923.75
$ws.Range("I80").Value = 1848.6666
$ws.Range("K80").Value = 1848.6666
$ws.Range("M80").Value = -850.6666

$ws.Range("H83").Value = 1923.75
$ws.Range("I83").Value = 1848.6666
$ws.Range("K83").Value = 9243.333000000001
$ws.Range("M83").Value = -4251.333000000001

$ws.Range("H92").Value = 28752.824
$ws.Range("J92").Value = 27260.666
$ws.Range("L92").Value = 27260.666
$ws.Range("N92").Value = -31004.666

$ws.Range("H95").Value = 27000
$ws.Range("J95").Value = 27000
$ws.Range("L95").Value = 27000
$ws.Range("N95").Value = -32492

$ws.Range("H97").Value = 1055.5
$ws.Range("I97").Value = 858.1429000000001
$ws.Range("J97").Value = 1252.8572
$ws.Range("K97").Value = 858.1429000000001
$ws.Range("L97").Value = 1252.8572
$ws.Range("M97").Value = -362.1429000000001
$ws.Range("N97").Value = -2244.8572

$ws.Range("H102").Value = 3545.923
$ws.Range("I102").Value = 3545.923
$ws.Range("K102").Value = 3545.923
$ws.Range("M102").Value = -1923.923

$ws.Range("H122").Value = 1379.5
$ws.Range("I122").Value = 1379.5
$ws.Range("K122").Value = 4138.5
$ws.Range("M122").Value = -1688.5

$ws.Range("H126").Value = 7848.852
$ws.Range("I126").Value = 7856
$ws.Range("J126").Value = 7828.4287
$ws.Range("K126").Value = 23568
$ws.Range("L126").Value = 23485.2861
$ws.Range("M126").Value = -21098
$ws.Range("N126").Value = -28425.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3682.2083
$ws.Range("I40").Value = 3899.0952
$ws.Range("K40").Value = 3899.0952
$ws.Range("M40").Value = -3763.0952

$ws.Range("H82").Value = 1762.3823
$ws.Range("I82").Value = 1456.8125
$ws.Range("J82").Value = 2034
$ws.Range("K82").Value = 1456.8125
$ws.Range("L82").Value = 2034
$ws.Range("M82").Value = -1095.8125
$ws.Range("N82").Value = -2756

$ws.Range("H85").Value = 1762.3823
$ws.Range("I85").Value = 1456.8125
$ws.Range("J85").Value = 2034
$ws.Range("K85").Value = 1456.8125
$ws.Range("L85").Value = 2034
$ws.Range("M85").Value = -208.8125
$ws.Range("N85").Value = -4530

$ws.Range("H122").Value = 4615.9614
$ws.Range("I122").Value = 3681.2354
$ws.Range("K122").Value = 11043.7062
$ws.Range("M122").Value = -8593.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 39847.375
$ws.Range("I61").Value = 39847.375
$ws.Range("K61").Value = 39847.375
$ws.Range("M61").Value = -39555.375

$ws.Range("H107").Value = 1407.4482
$ws.Range("I107").Value = 1143.5294
$ws.Range("J107").Value = 1781.3334
$ws.Range("K107").Value = 3430.5882
$ws.Range("L107").Value = 5344.0002
$ws.Range("M107").Value = -1510.5882
$ws.Range("N107").Value = -9184.0002

$ws.Range("H122").Value = 3108.9473
$ws.Range("I122").Value = 2515.6667
$ws.Range("J122").Value = 3382.7693
$ws.Range("K122").Value = 7547.000100000001
$ws.Range("L122").Value = 10148.3079
$ws.Range("M122").Value = -5097.000100000001
$ws.Range("N122").Value = -15048.3079

$ws.Range("H126").Value = 7106.231
$ws.Range("I126").Value = 4833.7173
$ws.Range("K126").Value = 14501.1519
$ws.Range("M126").Value = -12031.1519
